# js-bootcamp.pptx - "adding more notes and minor tweaks to presentation"
#
# 1) Slide 2 ("A brief history"): the bullet that used to read
#       "..., finally JavaScript"
#    becomes two sentences:
#       "..., finally names JavaScript when Sun became involved"
#
# 2) Slide 3 ("Some basic facts"): the two runs
#       "Though not without " + "some ugly bits"
#    are merged back into a single run of text
#       "Though not without some ugly bits"

$p = $ppt.ActivePresentation

function Get-ContentPlaceholder($slide) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq "Content Placeholder 2") {
            return $shp
        }
    }
    return $slide.Shapes.Item(2)
}

# ---------------------------------------------------------------------------
# Slide 2: split ", finally JavaScript" into two runs.
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = Get-ContentPlaceholder $slide2
$tr2 = $shape2.TextFrame.TextRange

$oldTail = ", finally JavaScript"
$newTail = ", finally names JavaScript when Sun "
$newRun = "became involved"

$full2 = $tr2.Text
$pos = $full2.IndexOf($oldTail)
if ($pos -ge 0) {
    # Rewrite the existing run's text in-place; this keeps its original
    # run properties (dirty="0" smtClean="0") untouched.
    $existingRun = $tr2.Characters($pos + 1, $oldTail.Length)
    $existingRun.Text = $newTail

    # Find which paragraph now holds that text so we can append a brand
    # new sibling run right after it (InsertAfter on a Paragraphs() range
    # creates a genuinely new <a:r>, instead of merging into the
    # preceding run the way Characters().InsertAfter does).
    $targetParaIndex = -1
    for ($i = 1; $i -le $tr2.Paragraphs().Count; $i++) {
        $para = $tr2.Paragraphs($i, 1)
        if ($para.Start -le ($pos + 1) -and ($para.Start + $para.Length) -gt $pos) {
            $targetParaIndex = $i
            break
        }
    }

    if ($targetParaIndex -gt 0) {
        $targetPara = $tr2.Paragraphs($targetParaIndex, 1)
        $targetPara.InsertAfter($newRun) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# Slide 3: merge "Though not without " + "some ugly bits" into one run.
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$shape3 = Get-ContentPlaceholder $slide3
$tr3 = $shape3.TextFrame.TextRange

$mergedOld = "Though not without " + "some ugly bits"
$mergedNew = "Though not without some ugly bits"

$full3 = $tr3.Text
$pos3 = $full3.IndexOf($mergedOld)
if ($pos3 -ge 0) {
    $mergeRange = $tr3.Characters($pos3 + 1, $mergedOld.Length)
    $mergeRange.Text = $mergedNew
}
